$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 90.860471916029525
$ws.Range("C2").Value = 85.278801938916132
$ws.Range("D2").Value = 48.813148317330196
$ws.Range("E2").Value = 73.305584044857653
$ws.Range("F2").Value = 66.216947413160355
$ws.Range("G2").Value = 67.229567205124212
$ws.Range("H2").Value = 91.374483818770145
$ws.Range("I2").Value = 81.534821816404374
$ws.Range("J2").Value = 76.685508077030491
$ws.Range("K2").Value = 79.93213317637381
$ws.Range("L2").Value = 85.137707370925668
$ws.Range("M2").Value = 68.83594907468779
$ws.Range("N2").Value = 83.685038584001106
$ws.Range("O2").Value = 55.699168451803679
$ws.Range("P2").Value = 101.27883778950699
$ws.Range("Q2").Value = 85.611492581789363
$ws.Range("R2").Value = 62.140063025524782
$ws.Range("T2").Value = 77.028511824947842
$ws.Range("U2").Value = 60.885628682385331
$ws.Range("V2").Value = 67.836789263764885
$ws.Range("W2").Value = 69.676829265588125
$ws.Range("X2").Value = 81.548860354004375
$ws.Range("Y2").Value = 55.67785536847164
$ws.Range("Z2").Value = 31.934854842981782
$ws.Range("AA2").Value = 42.75393681786538
$ws.Range("AB2").Value = 22.7815070975686
$ws.Range("AC2").Value = 52.993651277176603
$ws.Range("AD2").Value = 35.506720525962585
$ws.Range("AE2").Value = 66.585393887098888
$ws.Range("AF2").Value = 56.26454766936105
$ws.Range("AG2").Value = 58.598621544015707
$ws.Range("AH2").Value = 50.526662582709122
$ws.Range("AI2").Value = 30.751748302656452
$ws.Range("AJ2").Value = 58.563490696825696
$ws.Range("AK2").Value = 76.875347588072685
$ws.Range("AL2").Value = 79.77470128369977
$ws.Range("AM2").Value = 52.557941537573328
$ws.Range("AN2").Value = 51.337072060398604
$ws.Range("AO2").Value = 60.744154816114872
$ws.Range("AP2").Value = 57.329891261222045
$ws.Range("AQ2").Value = 82.065695254338223
$ws.Range("AR2").Value = 68.912156582675607
$ws.Range("AS2").Value = 74.551903875929852
$ws.Range("AT2").Value = 33.459558556051206
$ws.Range("AU2").Value = 73.63471513656971
$ws.Range("AV2").Value = 44.102543507086345
$ws.Range("AW2").Value = 41.173723408646588
$ws.Range("AX2").Value = 42.371492852737681
$ws.Range("AY2").Value = 44.613730176189321
$ws.Range("B3").Value = 83.543010585939058
$ws.Range("C3").Value = 70.217220856918715
$ws.Range("D3").Value = 57.357627020653325
$ws.Range("E3").Value = -58.467808491488135
$ws.Range("F3").Value = 75.56371466303068
$ws.Range("G3").Value = 66.922366093784944
$ws.Range("H3").Value = 75.705753328229733
$ws.Range("I3").Value = 88.687934623702503
$ws.Range("J3").Value = 69.680620262535669
$ws.Range("K3").Value = 73.581725707516028
$ws.Range("L3").Value = 74.048132536449302
$ws.Range("M3").Value = 71.074851363801471
$ws.Range("N3").Value = 72.869241919010591
$ws.Range("O3").Value = 31.908378129429536
$ws.Range("P3").Value = 49.437117323992105
$ws.Range("Q3").Value = 87.270004518605219
$ws.Range("R3").Value = 62.424865999607469
$ws.Range("S3").Value = 71.341814128989412
$ws.Range("T3").Value = 66.779782326304769
$ws.Range("U3").Value = 62.696918311791535
$ws.Range("V3").Value = 75.477605191000691
$ws.Range("W3").Value = 77.782815109556893
$ws.Range("X3").Value = 69.921751962487917
$ws.Range("Y3").Value = 30.982515943128309
$ws.Range("Z3").Value = 30.994613931163045
$ws.Range("AA3").Value = 50.449498435488607
$ws.Range("AB3").Value = 46.501664063050903
$ws.Range("AC3").Value = 65.628076201330501
$ws.Range("AD3").Value = 46.881974754576255
$ws.Range("AE3").Value = 53.617350632880175
$ws.Range("AF3").Value = 38.594192082374015
$ws.Range("AG3").Value = 68.279259910480178
$ws.Range("AH3").Value = 73.10096711937716
$ws.Range("AI3").Value = 45.38232623781763
$ws.Range("AJ3").Value = 44.322034929131945
$ws.Range("AK3").Value = 60.23220289210505
$ws.Range("AL3").Value = 54.071201147366068
$ws.Range("AM3").Value = 42.837996016753941
$ws.Range("AN3").Value = 54.722123176417512
$ws.Range("AO3").Value = 14.543276288253479
$ws.Range("AP3").Value = 50.845552338048748
$ws.Range("AQ3").Value = 82.332743271861645
$ws.Range("AR3").Value = 66.128511354835169
$ws.Range("AS3").Value = 57.845659311851684
$ws.Range("AT3").Value = 51.624687221853861
$ws.Range("AU3").Value = 80.110665283159804
$ws.Range("AV3").Value = 55.275102946755602
$ws.Range("AW3").Value = 61.84281797009281
$ws.Range("AX3").Value = 68.935136697682779
$ws.Range("AY3").Value = 33.26895842529256
